$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1452
$ws.Range("E2").Value = 120
$ws.Range("F2").Value = 128
$ws.Range("G2").Value = 109
$ws.Range("H2").Value = 91
$ws.Range("I2").Value = 88
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 1678
$ws.Range("L2").Value = 528
$ws.Range("M2").Value = 1150
$ws.Range("N2").Value = 1129
$ws.Range("O2").Value = 21
$ws.Range("P2").Value = 208
$ws.Range("Q2").Value = 117
$ws.Range("R2").Value = -133
$ws.Range("S2").Value = -11
$ws.Range("T2").Value = 9
$ws.Range("U2").Value = 108
$ws.Range("V2").Value = 271
$ws.Range("W2").Value = 8.24
$ws.Range("X2").Value = 6.28
$ws.Range("Y2").Value = 8.01
$ws.Range("Z2").Value = 5.5
$ws.Range("AA2").Value = 45.9
$ws.Range("AB2").Value = 477.04
$ws.Range("AC2").Value = 210
$ws.Range("AD2").Value = 17.14
$ws.Range("AE2").Value = 3034
$ws.Range("AF2").Value = 1.19
$ws.Range("AG2").Value = 50
$ws.Range("AH2").Value = 1.39
$ws.Range("AI2").Value = 21.26
$ws.Range("AJ2").Value = 41678175

# Row 3
$ws.Range("D3").Value = 1334
$ws.Range("E3").Value = 120
$ws.Range("F3").Value = 120
$ws.Range("G3").Value = 167
$ws.Range("H3").Value = 142
$ws.Range("I3").Value = 138
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 1755
$ws.Range("L3").Value = 458
$ws.Range("M3").Value = 1296
$ws.Range("N3").Value = 1275
$ws.Range("O3").Value = 21
$ws.Range("P3").Value = 208
$ws.Range("Q3").Value = 171
$ws.Range("R3").Value = -142
$ws.Range("S3").Value = -58
$ws.Range("T3").Value = 10
$ws.Range("U3").Value = 160
$ws.Range("V3").Value = 237
$ws.Range("W3").Value = 8.97
$ws.Range("X3").Value = 10.63
$ws.Range("Y3").Value = 11.5
$ws.Range("Z3").Value = 8.26
$ws.Range("AA3").Value = 35.36
$ws.Range("AB3").Value = 535.49
$ws.Range("AC3").Value = 332
$ws.Range("AD3").Value = 12.96
$ws.Range("AE3").Value = 3426
$ws.Range("AF3").Value = 1.26
$ws.Range("AG3").Value = 50
$ws.Range("AH3").Value = 1.16
$ws.Range("AI3").Value = 13.46
$ws.Range("AJ3").Value = 41678175

# Row 4
$ws.Range("D4").Value = 1337
$ws.Range("E4").Value = 154
$ws.Range("F4").Value = 162
$ws.Range("G4").Value = 172
$ws.Range("H4").Value = 125
$ws.Range("I4").Value = 120
$ws.Range("J4").Value = 6
$ws.Range("K4").Value = 1878
$ws.Range("L4").Value = 518
$ws.Range("M4").Value = 1360
$ws.Range("N4").Value = 1338
$ws.Range("O4").Value = 22
$ws.Range("P4").Value = 208
$ws.Range("Q4").Value = 92
$ws.Range("R4").Value = -87
$ws.Range("S4").Value = -2
$ws.Range("T4").Value = 9
$ws.Range("U4").Value = 83
$ws.Range("V4").Value = 254
$ws.Range("W4").Value = 11.5
$ws.Range("X4").Value = 9.38
$ws.Range("Y4").Value = 9.16
$ws.Range("Z4").Value = 6.9
$ws.Range("AA4").Value = 38.07
$ws.Range("AB4").Value = 589.97
$ws.Range("AC4").Value = 287
$ws.Range("AD4").Value = 12.81
$ws.Range("AE4").Value = 3596
$ws.Range("AF4").Value = 1.02
$ws.Range("AG4").Value = 70
$ws.Range("AH4").Value = 1.9
$ws.Range("AI4").Value = 21.76
$ws.Range("AJ4").Value = 41678175

# Row 5
$ws.Range("D5").Value = 1548
$ws.Range("E5").Value = 122
$ws.Range("F5").Value = 122
$ws.Range("G5").Value = 123
$ws.Range("H5").Value = 73
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 1900
$ws.Range("L5").Value = 513
$ws.Range("M5").Value = 1387
$ws.Range("N5").Value = 1366
$ws.Range("O5").Value = 21
$ws.Range("P5").Value = 208
$ws.Range("Q5").Value = -98
$ws.Range("R5").Value = 25
$ws.Range("S5").Value = 32
$ws.Range("T5").Value = 23
$ws.Range("U5").Value = -121
$ws.Range("V5").Value = 261
$ws.Range("W5").Value = 7.88
$ws.Range("X5").Value = 4.7
$ws.Range("Y5").Value = 5.17
$ws.Range("Z5").Value = 3.85
$ws.Range("AA5").Value = 37
$ws.Range("AB5").Value = 618.74
$ws.Range("AC5").Value = 168
$ws.Range("AD5").Value = 22.92
$ws.Range("AE5").Value = 3669
$ws.Range("AF5").Value = 1.05
$ws.Range("AG5").Value = 70
$ws.Range("AH5").Value = 1.82
$ws.Range("AI5").Value = 37.32
$ws.Range("AJ5").Value = 41678175

# Row 6
$ws.Range("D6").Value = 1577
$ws.Range("E6").Value = 133
$ws.Range("F6").Value = 133
$ws.Range("G6").Value = 152
$ws.Range("H6").Value = 116
$ws.Range("I6").Value = 113
$ws.Range("K6").Value = 1937
$ws.Range("L6").Value = 460
$ws.Range("M6").Value = 1476
$ws.Range("N6").Value = 1455
$ws.Range("P6").Value = 208
$ws.Range("Q6").Value = 92
$ws.Range("R6").Value = 65
$ws.Range("S6").Value = -81
$ws.Range("T6").Value = 10
$ws.Range("U6").Value = 82
$ws.Range("V6").Value = 209
$ws.Range("W6").Value = 8.44
$ws.Range("X6").Value = 7.36
$ws.Range("Y6").Value = 7.99
$ws.Range("Z6").Value = 6.05
$ws.Range("AA6").Value = 31.19
$ws.Range("AB6").Value = 659.67
$ws.Range("AC6").Value = 270
$ws.Range("AD6").Value = 21.11
$ws.Range("AE6").Value = 3908
$ws.Range("AF6").Value = 1.46
$ws.Range("AG6").Value = 70
$ws.Range("AH6").Value = 1.23
$ws.Range("AI6").Value = 23.11
$ws.Range("AJ6").Value = 41678175

# Rows 7-9 held estimate/forecast figures (2019E-2021E) that are no longer
# valid, so every data column is blanked out entirely, leaving only the
# row index (A) and the year/metric labels (B, C).
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
